$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 248, shifting existing rows 248:290 down to 249:291
$ws.Range("A248").EntireRow.Insert()

# Populate the newly inserted row 248 with the new weekly price record
$ws.Range("A248").Value = 4
$ws.Range("B248").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C248").Value = "Los Lagos"
$ws.Range("D248").Value = 44984
$ws.Range("E248").Value = 10
$ws.Range("F248").Value = "Fruta"
$ws.Range("G248").Value = 100103
$ws.Range("H248").Value = "Frutos de hueso (carozo)"
$ws.Range("I248").Value = 100103002
$ws.Range("J248").Value = "Ciruela"
$ws.Range("K248").Value = "Larry Ann"
$ws.Range("L248").Value = "Primera"
$ws.Range("M248").Value = 200
$ws.Range("N248").Value = 16000
$ws.Range("O248").Value = 17000
$ws.Range("P248").Value = 16500
$ws.Range("Q248").Value = "$/caja 14 kilos granel"
$ws.Range("R248").Value = "Región de O'Higgins"
$ws.Range("S248").Value = 1179
$ws.Range("T248").Value = 14
